$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New month of charges (2024-09-01, serial 45536) mirroring the layout of
# the previous month's block (rows 11-19), appended as rows 20-28.

$newRows = @(
    @{ Row = 20; A = "electricite";          B = 71  },
    @{ Row = 21; A = "gaz";                  B = 22  },
    @{ Row = 22; A = "copro";                B = 281 },
    @{ Row = 23; A = "box ";                 B = 30  },
    @{ Row = 24; A = "credit et assurances"; B = 400 },
    @{ Row = 25; A = "marceline";            B = 0   },
    @{ Row = 26; A = "charbel";              B = 0   },
    @{ Row = 27; A = "adekemi";              B = 0   },
    @{ Row = 28; A = "kadi";                 B = 0   }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $srcRow = $rowNum - 9

    # Copy formatting (number format, font, style indices) from the
    # corresponding row in the previous month's block so no new styles
    # are introduced.
    $ws.Range("A" + $srcRow + ":D" + $srcRow).Copy() | Out-Null
    $ws.Range("A" + $rowNum + ":D" + $rowNum).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = 45536
    $ws.Range("D" + $rowNum).Value = "Appartement Madoumier T4"
}

$excel.CutCopyMode = 0

$ws.Range("A1").Select() | Out-Null
$ws.Range("I30").Select() | Out-Null
